# Scheduled market-data refresh: update computed profit columns (H:N) for the
# affected leve rows across all eight crafting-job sheets. Values below are the
# latest pulled averages/prices; no formulas are involved (source data is static).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2687.625
$ws.Range("I32").Value = 1533.6666
$ws.Range("J32").Value = 3380
$ws.Range("K32").Value = 1533.6666
$ws.Range("L32").Value = 3380
$ws.Range("M32").Value = -1207.6666
$ws.Range("N32").Value = -4032
$ws.Range("H40").Value = 2225.3157
$ws.Range("I40").Value = 1820.25
$ws.Range("J40").Value = 2333.3333
$ws.Range("K40").Value = 1820.25
$ws.Range("L40").Value = 2333.3333
$ws.Range("M40").Value = -1645.25
$ws.Range("N40").Value = -2683.3333
$ws.Range("H116").Value = 2845.3635
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 2659.8
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 2659.8
$ws.Range("M116").Value = 442
$ws.Range("N116").Value = -9543.799999999999
$ws.Range("H125").Value = 6762.5625
$ws.Range("I125").Value = 2939.0908
$ws.Range("J125").Value = 15174.2
$ws.Range("K125").Value = 26451.8172
$ws.Range("L125").Value = 136567.8
$ws.Range("M125").Value = -23991.8172
$ws.Range("N125").Value = -141487.8
$ws.Range("H129").Value = 999.0238000000001
$ws.Range("J129").Value = 1153.8182
$ws.Range("L129").Value = 3461.4546
$ws.Range("N129").Value = -13461.4546
$ws.Range("H132").Value = 2908.457
$ws.Range("I132").Value = 2955.5557
$ws.Range("J132").Value = 2749.5
$ws.Range("K132").Value = 8866.667099999999
$ws.Range("L132").Value = 8248.5
$ws.Range("M132").Value = -6336.667099999999
$ws.Range("N132").Value = -13308.5
$ws.Range("H138").Value = 3113.19
$ws.Range("I138").Value = 767.6
$ws.Range("J138").Value = 3895.0532
$ws.Range("K138").Value = 2302.8
$ws.Range("L138").Value = 11685.1596
$ws.Range("M138").Value = 2837.2
$ws.Range("N138").Value = -21965.1596

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 824
$ws.Range("I61").Value = 706.4
$ws.Range("K61").Value = 706.4
$ws.Range("M61").Value = -494.4
$ws.Range("H74").Value = 1142.8334
$ws.Range("I74").Value = 1106.3334
$ws.Range("J74").Value = 1252.3334
$ws.Range("K74").Value = 1106.3334
$ws.Range("L74").Value = 1252.3334
$ws.Range("M74").Value = -232.3334
$ws.Range("N74").Value = -3000.3334
$ws.Range("H77").Value = 1142.8334
$ws.Range("I77").Value = 1106.3334
$ws.Range("J77").Value = 1252.3334
$ws.Range("K77").Value = 5531.666999999999
$ws.Range("L77").Value = 6261.666999999999
$ws.Range("M77").Value = -1163.666999999999
$ws.Range("N77").Value = -14997.667
$ws.Range("H122").Value = 951.7143
$ws.Range("I122").Value = 943.6667
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2831.0001
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -381.0001000000002
$ws.Range("N122").Value = -7900
$ws.Range("H136").Value = 824
$ws.Range("I136").Value = 706.4
$ws.Range("K136").Value = 2119.2
$ws.Range("M136").Value = 430.8000000000002

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 429.6087
$ws.Range("I80").Value = 248.875
$ws.Range("J80").Value = 526
$ws.Range("K80").Value = 248.875
$ws.Range("L80").Value = 526
$ws.Range("M80").Value = 749.125
$ws.Range("N80").Value = -2522
$ws.Range("H83").Value = 429.6087
$ws.Range("I83").Value = 248.875
$ws.Range("J83").Value = 526
$ws.Range("K83").Value = 1244.375
$ws.Range("L83").Value = 2630
$ws.Range("M83").Value = 3747.625
$ws.Range("N83").Value = -12614
$ws.Range("H99").Value = 2155.3333
$ws.Range("I99").Value = 2060
$ws.Range("J99").Value = 2386.8572
$ws.Range("K99").Value = 2060
$ws.Range("L99").Value = 2386.8572
$ws.Range("M99").Value = -562
$ws.Range("N99").Value = -5382.8572
$ws.Range("H140").Value = 72424.75
$ws.Range("J140").Value = 72424.75
$ws.Range("L140").Value = 72424.75
$ws.Range("N140").Value = -82784.75

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3187.5312
$ws.Range("J31").Value = 2638
$ws.Range("L31").Value = 2638
$ws.Range("N31").Value = -3228
$ws.Range("H34").Value = 3187.5312
$ws.Range("J34").Value = 2638
$ws.Range("L34").Value = 2638
$ws.Range("N34").Value = -3042
$ws.Range("H105").Value = 1995
$ws.Range("I105").Value = 2162.2222
$ws.Range("J105").Value = 1493.3334
$ws.Range("K105").Value = 2162.2222
$ws.Range("L105").Value = 1493.3334
$ws.Range("M105").Value = -415.2222000000002
$ws.Range("N105").Value = -4987.3334
$ws.Range("H132").Value = 1556.8286
$ws.Range("I132").Value = 813.875
$ws.Range("J132").Value = 3177.818
$ws.Range("K132").Value = 2441.625
$ws.Range("L132").Value = 9533.454000000002
$ws.Range("M132").Value = 88.375
$ws.Range("N132").Value = -14593.454
$ws.Range("H134").Value = 2190.2632
$ws.Range("I134").Value = 1661
$ws.Range("K134").Value = 4983
$ws.Range("M134").Value = -2448

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6226.364
$ws.Range("J80").Value = 6226.364
$ws.Range("L80").Value = 18679.092
$ws.Range("N80").Value = -20551.092
$ws.Range("H83").Value = 6226.364
$ws.Range("J83").Value = 6226.364
$ws.Range("L83").Value = 56037.276
$ws.Range("N83").Value = -65397.276
$ws.Range("H113").Value = 665.5714
$ws.Range("J113").Value = 675.71155
$ws.Range("L113").Value = 2027.13465
$ws.Range("N113").Value = -6367.13465
$ws.Range("H129").Value = 30524
$ws.Range("J129").Value = 60202.234
$ws.Range("L129").Value = 180606.702
$ws.Range("N129").Value = -190606.702

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6010.75
$ws.Range("I70").Value = 5017.2
$ws.Range("K70").Value = 5017.2
$ws.Range("M70").Value = -4747.2
$ws.Range("H73").Value = 6010.75
$ws.Range("I73").Value = 5017.2
$ws.Range("K73").Value = 5017.2
$ws.Range("M73").Value = -4081.2
$ws.Range("H126").Value = 2505.1
$ws.Range("I126").Value = 2191.25
$ws.Range("K126").Value = 6573.75
$ws.Range("M126").Value = -4103.75
$ws.Range("H132").Value = 3189.7693
$ws.Range("I132").Value = 3008
$ws.Range("J132").Value = 3533.111
$ws.Range("K132").Value = 9024
$ws.Range("L132").Value = 10599.333
$ws.Range("M132").Value = -6494
$ws.Range("N132").Value = -15659.333

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1424.5
$ws.Range("I22").Value = 1250.25
$ws.Range("J22").Value = 1511.625
$ws.Range("K22").Value = 1250.25
$ws.Range("L22").Value = 1511.625
$ws.Range("M22").Value = -955.25
$ws.Range("N22").Value = -2101.625
$ws.Range("H27").Value = 1424.5
$ws.Range("I27").Value = 1250.25
$ws.Range("J27").Value = 1511.625
$ws.Range("K27").Value = 1250.25
$ws.Range("L27").Value = 1511.625
$ws.Range("M27").Value = -1143.25
$ws.Range("N27").Value = -1725.625
$ws.Range("H92").Value = 18784.143
$ws.Range("J92").Value = 18784.143
$ws.Range("L92").Value = 18784.143
$ws.Range("N92").Value = -23776.143
$ws.Range("H136").Value = 4684.2085
$ws.Range("I136").Value = 1521.05
$ws.Range("J136").Value = 20500
$ws.Range("K136").Value = 4563.15
$ws.Range("L136").Value = 61500
$ws.Range("M136").Value = -2013.15
$ws.Range("N136").Value = -66600

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 40027.5
$ws.Range("I34").Value = 10026
$ws.Range("J34").Value = 70029
$ws.Range("K34").Value = 10026
$ws.Range("L34").Value = 70029
$ws.Range("M34").Value = -9823
$ws.Range("N34").Value = -70435
$ws.Range("H37").Value = 70029
$ws.Range("J37").Value = 70029
$ws.Range("L37").Value = 70029
$ws.Range("N37").Value = -70435
$ws.Range("H40").Value = 3500.5
$ws.Range("I40").Value = 2001
$ws.Range("K40").Value = 2001
$ws.Range("M40").Value = -1852
$ws.Range("H87").Value = 18000
$ws.Range("J87").Value = 18000
$ws.Range("L87").Value = 18000
$ws.Range("N87").Value = -20496
$ws.Range("H90").Value = 18000
$ws.Range("J90").Value = 18000
$ws.Range("L90").Value = 54000
$ws.Range("N90").Value = -66480
$ws.Range("H132").Value = 1061.9375
$ws.Range("I132").Value = 1043.2222
$ws.Range("J132").Value = 1163
$ws.Range("K132").Value = 3129.6666
$ws.Range("L132").Value = 3489
$ws.Range("M132").Value = -599.6665999999996
$ws.Range("N132").Value = -8549
